# Inserts one new data row at row 514 on the (single) worksheet of the
# "Hortaliza, Macroferia Regional de Talca - Brócoli" workbook.
#
# Effect: every existing data row from 514 down to 592 is pushed down by
# one row (514->515, 515->516, ..., 592->593, exactly like Excel's native
# "Insert Row" behaviour), and the vacated row 514 is populated with a
# brand-new record (date serial 45127, i.e. 2023-07-20).
#
# This grows the used range from A1:R592 to A1:R593.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 514..592 down to 515..593, leaving a blank row 514 behind
# (row 514 inherits formatting from the row above it, same as Excel).
$ws.Rows.Item(514).Insert()

# Populate the newly inserted row 514 with the new record.
$ws.Cells.Item(514, 1).Value  = 5
$ws.Cells.Item(514, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(514, 3).Value  = "Maule"
$ws.Cells.Item(514, 4).Value  = 45127
$ws.Cells.Item(514, 5).Value  = 7
$ws.Cells.Item(514, 6).Value  = 100112023
$ws.Cells.Item(514, 7).Value  = "Brócoli"
$ws.Cells.Item(514, 8).Value  = "Sin especificar"
$ws.Cells.Item(514, 9).Value  = "Primera"
$ws.Cells.Item(514, 10).Value = 3000
$ws.Cells.Item(514, 11).Value = 600
$ws.Cells.Item(514, 12).Value = 600
$ws.Cells.Item(514, 13).Value = 600
$ws.Cells.Item(514, 14).Value = "`$/unidad"
$ws.Cells.Item(514, 15).Value = "Región del Maule"
$ws.Cells.Item(514, 16).Value = 600
$ws.Cells.Item(514, 17).Value = 1
$ws.Cells.Item(514, 18).Value = "Hortaliza"
